$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prime rows 13-16 so their row "spans" memory extends to column D,
#     matching the source workbook (those rows briefly touched D before
#     being cleared again). Do this before filling real data.
$ws.Range("D13").Value = "x"
$ws.Range("D13").Value = ""
$ws.Range("D14").Value = "x"
$ws.Range("D14").Value = ""
$ws.Range("D15").Value = "x"
$ws.Range("D15").Value = ""
$ws.Range("D16").Value = "x"
$ws.Range("D16").Value = ""

# --- Column A (dates), rows 13-26, top to bottom. Rows 19 & 21 are true
#     numeric Excel dates (style copied from A2 so they reuse numFmtId 14);
#     the rest are literal text that happens to look like a date.
$ws.Range("A13").Value = "18/12/2019"
$ws.Range("A14").Value = "19/12/2019"
$ws.Range("A15").Value = "20/12/2019"
$ws.Range("A16").Value = "21/12/2019"
$ws.Range("A17").Value = "22/12/2019"
$ws.Range("A18").Value = "23/12/2019"

$ws.Range("A2").Copy($ws.Range("A19"))
$ws.Range("A19").Value = 43862

$ws.Range("A20").Value = " 6/1/2020"

$ws.Range("A2").Copy($ws.Range("A21"))
$ws.Range("A21").Value = 44013

$ws.Range("A22").Value = " 8/1/2020"
$ws.Range("A23").Value = " 9/1/2020"
$ws.Range("A24").Value = " 10/1/2020"
$ws.Range("A25").Value = "20/1/2020"
$ws.Range("A26").Value = "22/1/2020"

# --- Column B (task text), rows 13-26, top to bottom.
$ws.Range("B13").Value = "design the home page"
$ws.Range("B14").Value = "edit the document"
$ws.Range("B15").Value = "full stack management"
$ws.Range("B16").Value = "design the home page"
$ws.Range("B17").Value = "design ui"
$ws.Range("B18").Value = "design"
$ws.Range("B19").Value = "prepare for first review "
$ws.Range("B20").Value = "basic java concepts"
$ws.Range("B21").Value = "java concepts"
$ws.Range("B22").Value = "design the front page "
$ws.Range("B23").Value = "design the navigation"
$ws.Range("B24").Value = "design the test concept"
$ws.Range("B25").Value = "design the test concept"
$ws.Range("B26").Value = "navigation"

# --- Column A (dates), rows 27-33, top to bottom (all literal text).
$ws.Range("A27").Value = "23/1/2020"
$ws.Range("A28").Value = "24/1/2020"
$ws.Range("A29").Value = "27/1/2020"
$ws.Range("A30").Value = "28/1/2020"
$ws.Range("A31").Value = "29/1/2020"
$ws.Range("A32").Value = "30/1/2020"
$ws.Range("A33").Value = "31/1/2020"

# --- Column B (task text), rows 27-37, top to bottom.
$ws.Range("B27").Value = "DAO"
$ws.Range("B28").Value = "corrected errors in DAO"
$ws.Range("B29").Value = "Did dao implementation"
$ws.Range("B30").Value = "controller"
$ws.Range("B31").Value = "main controller"
$ws.Range("B32").Value = "created mode"
$ws.Range("B33").Value = "database"
$ws.Range("B34").Value = "insert statement"
$ws.Range("B35").Value = "jdbc"
$ws.Range("B36").Value = "localhost run"
$ws.Range("B37").Value = "corrected errors in project"

# --- Column A (true numeric dates), rows 34-37 - style copied from A2 so
#     they reuse the existing date style (numFmtId 14) like A19/A21 above.
$ws.Range("A2").Copy($ws.Range("A34"))
$ws.Range("A34").Value = 43892

$ws.Range("A2").Copy($ws.Range("A35"))
$ws.Range("A35").Value = 43923

$ws.Range("A2").Copy($ws.Range("A36"))
$ws.Range("A36").Value = 43953

$ws.Range("A2").Copy($ws.Range("A37"))
$ws.Range("A37").Value = 43984

# --- View settings to match the target sheetView.
$ws.Application.ActiveWindow.Zoom = 85
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B37").Select()
